# Lesson3:Task1 Matrices, Tables - Multidimensional arrays
#
# Splits a handful of single runs into two runs (to mirror the edits the
# original author made while capitalising the first letter of some terms),
# capitalises a couple of first letters, and fills in the previously-empty
# 3rd data row of the dictionary table with "Multidimensional arrays" /
# "Многомерные массивы".

$d = $word.ActiveDocument

function Split-FirstChar($range1) {
    # Forces Word to break the run that $range1 (a 1-character range) lives
    # in into its own run, without altering the effective/visible
    # character formatting (the property is toggled back off immediately).
    $range1.Font.Bold = 1
    $range1.Font.Bold = 0
}

# --- Row 1: "bounds of array" -> "B" + "ounds of array" -------------------
$rng = $d.Content.Duplicate
$rng.Find.Execute("bounds of array", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$firstChar = $d.Range($rng.Start, $rng.Start + 1)
$firstChar.Text = "B"
Split-FirstChar $firstChar

# --- Row 1: "границы " -> "Г" + "раницы " ----------------------------------
$rng = $d.Content.Duplicate
$rng.Find.Execute("границы ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$firstChar = $d.Range($rng.Start, $rng.Start + 1)
$firstChar.Text = "Г"
Split-FirstChar $firstChar

# --- Row 2: "s" -> "S" (already its own run) -------------------------------
$rng = $d.Content.Duplicate
$rng.Find.Execute("square", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$firstChar = $d.Range($rng.Start, $rng.Start + 1)
$firstChar.Text = "S"

# --- Row 2: "квадрат" -> "К" + "вадрат" -------------------------------------
$rng = $d.Content.Duplicate
$rng.Find.Execute("квадрат", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$firstChar = $d.Range($rng.Start, $rng.Start + 1)
$firstChar.Text = "К"
Split-FirstChar $firstChar

# --- Row 3: fill in previously empty cells ---------------------------------
$table = $d.Tables.Item(1)

$cellEn = $table.Cell(3, 1)
$rngEn = $cellEn.Range
$rngEn.Collapse(0)
$rngEn.InsertBefore("Multidimensional arrays")

$cellRu = $table.Cell(3, 2)
$rngRu = $cellRu.Range
$rngRu.Collapse(0)
$rngRu.InsertBefore("Многомерные массивы")
